# Update column G ("K") values on the active worksheet to reflect the
# regenerated strike-count data (using K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G
$newValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 3
    13 = 0
    14 = 0
    15 = 2
    16 = 2
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 2
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 2
    29 = 4
    30 = 1
    31 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
